$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new command row (row 38): "sudo fdisk -l" / "vypis diskov v systeme"
$ws.Range("A38").Value = "sudo fdisk -l"
$ws.Range("B38").Value = "vypis diskov v systeme"

# Match the author's final selection / active cell
$ws.Range("B38").Select()
